# Fill in the trailing empty paragraph with the "Version management" text,
# building it up as separate runs the way Word naturally leaves them when
# you type a paragraph, split it, and then stitch the pieces back together
# by deleting the paragraph marks in between (Word does not re-merge runs
# across such a join).

function Join-WithNextParagraph($doc) {
    $count = $doc.Paragraphs.Count
    $prior = $doc.Paragraphs.Item($count - 1)
    $priorEnd = $prior.Range.End
    $mark = $doc.Range($priorEnd - 1, $priorEnd)
    $mark.Delete()
}

function Add-RunToLastParagraph($doc, $text) {
    $docEnd = $doc.Content.End
    $tail = $doc.Range($docEnd - 1, $docEnd - 1)
    $tail.InsertParagraphAfter()
    $doc.Paragraphs.Last.Range.Text = $text
    Join-WithNextParagraph $doc
}

$d = $word.ActiveDocument

$target = $d.Paragraphs.Last
$target.Range.Text = "When numerous developers are working on the same project, version control becomes a must."

Add-RunToLastParagraph $d " "
Add-RunToLastParagraph $d "Each developer may be making independent progress on their own version of a component at the same time."
Add-RunToLastParagraph $d " "
Add-RunToLastParagraph $d "In addition, previous versions are crucial in case the system has to roll back to them because of unexpected faults or misinterpreted modification requirements."

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
Write-Output "Final text: $($d.Paragraphs.Last.Range.Text)"
